$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

# Helper: columns that must contain an empty text value (matching the
# empty inlineStr cells already used throughout the sheet) rather than a
# truly blank/empty cell. Assigning a single leading apostrophe makes
# Excel store an empty text cell; resetting the style afterwards clears
# the transient "quote prefix" formatting flag so the cell matches the
# plain unstyled look of its neighbours.
$emptyTextCols = 2,9,10,11,12,13,17,19,23

foreach ($col in $emptyTextCols) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

$ws.Cells.Item($row, 1).Value = "2024-09-08 20:30:11"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 14).Value = 10
$ws.Cells.Item($row, 15).Value = 4
$ws.Cells.Item($row, 16).Value = 3
$ws.Cells.Item($row, 18).Value = 3
$ws.Cells.Item($row, 20).Value = 50
$ws.Cells.Item($row, 21).Value = 0
$ws.Cells.Item($row, 22).Value = "D:\Repositorio\jonatha1992\Predictor_ruleta\Data\Crupier.xlsx"
$ws.Cells.Item($row, 24).Value = "No es Simulación"
$ws.Cells.Item($row, 25).Value = 0
